$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 11 ("R40" rule) moves its label cell (B11) from the text "R40" to
# the text "1". A leading apostrophe forces Excel to store it as text
# (matching the shared-string / text cell the workbook ends up with)
# rather than silently converting it to the number 1.
$ws.Range("B11").Value = "'1"
